$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose Target cluster is "Resolving-Mac" (old rows 5 and 9).
# After deleting row 5, the old row 9 becomes row 8.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(8).Delete()

# Refresh the recomputed TPM-derived statistics for the remaining rows.
$ws.Range("M2").Value = 0.274713
$ws.Range("N2").Value = 0.824139
$ws.Range("O2").Value = 0.1055967877339779
$ws.Range("P2").Value = 0.1055967877339779
$ws.Range("Q2").Value = 0.043926700271
$ws.Range("R2").Value = 0.395340302439
$ws.Range("S2").Value = 0.00279596211882344
$ws.Range("T2").Value = 0.002795962118823441
$ws.Range("O3").Value = 0.3415960415058637
$ws.Range("P3").Value = 0.3415960415058638
$ws.Range("S3").Value = 0.009044684146988642
$ws.Range("T3").Value = 0.009044684146988646
$ws.Range("M4").Value = 1.438143333333333
$ws.Range("N4").Value = 4.31443
$ws.Range("O4").Value = 0.5528071707601584
$ws.Range("P4").Value = 0.5528071707601584
$ws.Range("Q4").Value = 0.2299595983811111
$ws.Range("R4").Value = 2.06963638543
$ws.Range("S4").Value = 0.01463707316886522
$ws.Range("T4").Value = 0.01463707316886522
$ws.Range("M5").Value = 0.274713
$ws.Range("N5").Value = 0.824139
$ws.Range("O5").Value = 0.1055967877339779
$ws.Range("P5").Value = 0.1055967877339779
$ws.Range("Q5").Value = 1.615079483376
$ws.Range("R5").Value = 14.535715350384
$ws.Range("S5").Value = 0.1028008256151544
$ws.Range("T5").Value = 0.1028008256151544
$ws.Range("O6").Value = 0.3415960415058637
$ws.Range("P6").Value = 0.3415960415058638
$ws.Range("S6").Value = 0.332551357358875
$ws.Range("T6").Value = 0.3325513573588751
$ws.Range("M7").Value = 1.438143333333333
$ws.Range("N7").Value = 4.31443
$ws.Range("O7").Value = 0.5528071707601584
$ws.Range("P7").Value = 0.5528071707601584
$ws.Range("Q7").Value = 8.455063254453334
$ws.Range("R7").Value = 76.09556929007999
$ws.Range("S7").Value = 0.5381700975912932
$ws.Range("T7").Value = 0.5381700975912932
